# Fixed update to excel issue
# 1) Rename the "Requested quantity" headers on the existing sheets.
# 2) Add a new "PO Forecast" sheet with a forecast table (ds, PO_Forecast,
#    yhat_lower, yhat_upper) formatted to match the other sheets.

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- Rename header labels -------------------------------------------------
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet at the end ----------------------------
$wsForecast = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match header formatting (bold, bordered, centered) used on the other sheets
$wsForecast.Range("A1:D1").Font.Bold          = $wsWeekly.Range("B1").Font.Bold
$wsForecast.Range("A1:D1").Borders.LineStyle  = $wsWeekly.Range("B1").Borders.LineStyle
$wsForecast.Range("A1:D1").HorizontalAlignment = $wsWeekly.Range("B1").HorizontalAlignment
$wsForecast.Range("A1:D1").VerticalAlignment   = $wsWeekly.Range("B1").VerticalAlignment

# Data rows
$data = @(
    @(45102.99999999999, 31, -33.72568847982005, 91.78819892675058),
    @(45144.99999999999, 42, -18.29927542506577, 98.17490041825384),
    @(45151.99999999999, 43, -18.28839434017091, 105.448901172796),
    @(45179.99999999999, 50, -9.876356475407478, 109.9233371697414),
    @(45186.99999999999, 52, -11.00060755492564, 114.7663868016302),
    @(45193.99999999999, 54, -4.157637551221536, 111.8315314658845),
    @(45200.99999999999, 55, -6.84890431044613, 109.5930757926303),
    @(45207.99999999999, 57, -1.782653178280052, 116.0403150146933),
    @(45221.99999999999, 60, -0.9368286267195063, 127.4525341274122),
    @(45228.99999999999, 62, 2.61283273288394, 127.4638495546624),
    @(45235.99999999999, 64, 10.06190116949698, 121.6630534874693),
    @(45242.99999999999, 66, 6.428242243603272, 127.8646598241972),
    @(45249.99999999999, 67, 9.280333330434525, 124.6319313497333),
    @(45263.99999999999, 71, 6.342119511132069, 124.1685540050262),
    @(45270.99999999999, 72, 5.394519175835359, 133.6459810312749),
    @(45277.99999999999, 74, 17.32790430232646, 136.6538386810329),
    @(45298.99999999999, 79, 15.65375005647252, 139.6462942445944),
    @(45305.99999999999, 81, 20.31825890469083, 140.3469720722117),
    @(45312.99999999999, 83, 20.61721830187597, 141.2963984722072),
    @(45319.99999999999, 84, 23.32895434327022, 146.5425036034977),
    @(45326.99999999999, 86, 27.55833376202148, 146.1768463350448),
    @(45333.99999999999, 88, 26.43188412858372, 148.5868962296653),
    @(45340.99999999999, 89, 26.17294499526075, 153.2523514159913),
    @(45347.99999999999, 91, 29.77213672832321, 152.1366440481529),
    @(45354.99999999999, 93, 29.54445523531914, 150.6017369074058)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Match the date-time number format used for the "ds" columns elsewhere
$wsForecast.Range("A2:A26").NumberFormat = $wsWeekly.Range("A2").NumberFormat

# Restore the originally active sheet/selection
$wsWeekly.Activate() | Out-Null
$wsWeekly.Range("A1").Select() | Out-Null
